$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; existing rows 11-70 shift down to 12-71
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the new weekly record
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"
$ws.Range("D11").Value = 44635
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100104
$ws.Range("H11").Value = "Frutos de pepita"
$ws.Range("I11").Value = 100104003
$ws.Range("J11").Value = "Membrillo"
$ws.Range("K11").Value = "Champion"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 400
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20500
$ws.Range("Q11").Value = "$/caja 18 kilos granel"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 1139
$ws.Range("T11").Value = 18
